
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Paridhi")

# --- Week 1 details (row 2, A2:J2 already merged) -------------------------
$ws.Range("A2:J2").Value = "We spent this week deciding our project. Our first idea was a Roomba (automated vaccuum and mop) so I spent this week reading about SLAM before deciding that it is out of the scope of this project for now. Then when we thought about doing Just Dance I read a little about the different ways we would track the dancer (imu vs camera) and how much memory we would need."
$ws.Range("A2:J2").WrapText = $true
$ws.Rows.Item(2).RowHeight = 49.2

$ws.Rows.Item(3).RowHeight = 15.6
$ws.Rows.Item(4).RowHeight = 15.6

# --- Week 2 details (row 6, A6:J6 already merged) --------------------------
$ws.Range("A6:J6").Value = "After we decided that we want to make Just Dance for our project, I read about the different ways to track body movements, image processing algorithms for detecting a body or parts of a body, and algorithms for comparing two videos. At the end of this week we had finalised that we wanted to do DTW for the comparision algorithm."
$ws.Range("A6:J6").WrapText = $true
$ws.Rows.Item(6).RowHeight = 49.2

# --- Remove the extra blank placeholder row before "Week 2" label ----------
$ws.Rows.Item(9).Delete()

# After the delete: row 9 = "Week 2" label, row 10 = the merged A10:J10 cell
# that used to hold the blank placeholder just under "Week 2".
$ws.Range("A10:J10").Value = "I read about how DTW works, and the different ways to make it faster, and smaller so that it can fit on an FPGA and work in real time. Decided on an implementation we liked."
$ws.Range("A10:J10").WrapText = $true
$ws.Rows.Item(10).RowHeight = 31.2

# --- Unmerge row 13 (was blank placeholder) and add new weekly notes -------
$ws.Range("A13:J13").UnMerge()
$ws.Range("A13").Value = "Week 3"
$ws.Range("A14").Value = "Implemented and tested a basic DTW and custom shift register. Will need to experiment with types of ditance and minimum units."

$ws.Range("A17").Value = "Weeks 4 and 5"
$ws.Range("A18").Value = "While testing DTW further, I realised that the basic version was missing half the functionality. Rewrote a more robust DTW and tested it. The score is much more accurate now."
$ws.Range("A19").Value = "Started working on the game loop."

# --- Sheet / selection state -------------------------------------------
$ws.Activate()
$ws.Range("A19").Select() | Out-Null
